$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)
$shape = $s.Shapes.Item(2)
$tr = $shape.TextFrame.TextRange
$para = $tr.Paragraphs(6)
$para.InsertAfter("s")
